$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 37039810
$ws.Range("I100").Value = 55557210
$ws.Range("K100").Value = 55557210
$ws.Range("M100").Value = -55556669
$ws.Range("H129").Value = 871.8795
$ws.Range("J129").Value = 872.9877
$ws.Range("L129").Value = 2618.9631
$ws.Range("N129").Value = -12618.9631
$ws.Range("H132").Value = 174164.83
$ws.Range("I132").Value = 174164.83
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 522494.49
$ws.Range("L132").Value = 0
$ws.Range("H138").Value = 2476.7334
$ws.Range("I138").Value = 1939
$ws.Range("J138").Value = 3014.4666
$ws.Range("K138").Value = 5817
$ws.Range("L138").Value = 9043.399800000001
$ws.Range("M138").Value = -677
$ws.Range("N138").Value = -19323.3998
$ws.Range("M132").Value = -519964.49

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6056.8
$ws.Range("I32").Value = 4016.205
$ws.Range("J32").Value = 13291.637
$ws.Range("K32").Value = 4016.205
$ws.Range("L32").Value = 13291.637
$ws.Range("M32").Value = -3729.205
$ws.Range("N32").Value = -13865.637
$ws.Range("H61").Value = 5819.067
$ws.Range("I61").Value = 4020.923
$ws.Range("J61").Value = 17507
$ws.Range("K61").Value = 4020.923
$ws.Range("L61").Value = 17507
$ws.Range("M61").Value = -3808.923
$ws.Range("N61").Value = -17931
$ws.Range("H74").Value = 1562.6471
$ws.Range("I74").Value = 1466
$ws.Range("J74").Value = 1876.75
$ws.Range("K74").Value = 1466
$ws.Range("L74").Value = 1876.75
$ws.Range("M74").Value = -592
$ws.Range("N74").Value = -3624.75
$ws.Range("H77").Value = 1562.6471
$ws.Range("I77").Value = 1466
$ws.Range("J77").Value = 1876.75
$ws.Range("K77").Value = 7330
$ws.Range("L77").Value = 9383.75
$ws.Range("M77").Value = -2962
$ws.Range("N77").Value = -18119.75
$ws.Range("H132").Value = 3544.6667
$ws.Range("I132").Value = 1595.091
$ws.Range("J132").Value = 5194.3076
$ws.Range("K132").Value = 4785.272999999999
$ws.Range("L132").Value = 15582.9228
$ws.Range("M132").Value = -2255.272999999999
$ws.Range("N132").Value = -20642.9228
$ws.Range("H136").Value = 5819.067
$ws.Range("I136").Value = 4020.923
$ws.Range("J136").Value = 17507
$ws.Range("K136").Value = 12062.769
$ws.Range("L136").Value = 52521
$ws.Range("M136").Value = -9512.769
$ws.Range("N136").Value = -57621

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3320.9524
$ws.Range("I134").Value = 1785.5714
$ws.Range("J134").Value = 4088.6428
$ws.Range("K134").Value = 5356.7142
$ws.Range("L134").Value = 12265.9284
$ws.Range("M134").Value = -2821.7142
$ws.Range("N134").Value = -17335.9284

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 795.8333
$ws.Range("I7").Value = 2104
$ws.Range("J7").Value = 292.69232
$ws.Range("K7").Value = 2104
$ws.Range("L7").Value = 292.69232
$ws.Range("M7").Value = -1991
$ws.Range("N7").Value = -518.69232
$ws.Range("H132").Value = 2854.3333
$ws.Range("I132").Value = 1677
$ws.Range("J132").Value = 3443
$ws.Range("K132").Value = 5031
$ws.Range("L132").Value = 10329
$ws.Range("M132").Value = -2501
$ws.Range("N132").Value = -15389

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1141.6061
$ws.Range("I5").Value = 671.8125
$ws.Range("J5").Value = 1583.7646
$ws.Range("K5").Value = 2015.4375
$ws.Range("L5").Value = 4751.293799999999
$ws.Range("M5").Value = -1903.4375
$ws.Range("N5").Value = -4975.293799999999
$ws.Range("H113").Value = 1645.1177
$ws.Range("I113").Value = 1607.4445
$ws.Range("J113").Value = 1687.5
$ws.Range("K113").Value = 4822.333500000001
$ws.Range("L113").Value = 5062.5
$ws.Range("M113").Value = -2652.333500000001
$ws.Range("N113").Value = -9402.5
$ws.Range("H122").Value = 728.8182
$ws.Range("I122").Value = 465
$ws.Range("J122").Value = 827.75
$ws.Range("K122").Value = 4185
$ws.Range("L122").Value = 7449.75
$ws.Range("M122").Value = -1735
$ws.Range("N122").Value = -12349.75
$ws.Range("H132").Value = 1014.94116
$ws.Range("I132").Value = 798.75
$ws.Range("J132").Value = 1207.1111
$ws.Range("K132").Value = 7188.75
$ws.Range("L132").Value = 10863.9999
$ws.Range("M132").Value = -4658.75
$ws.Range("N132").Value = -15923.9999
$ws.Range("H135").Value = 1141.6061
$ws.Range("I135").Value = 671.8125
$ws.Range("J135").Value = 1583.7646
$ws.Range("K135").Value = 6046.3125
$ws.Range("L135").Value = 14253.8814
$ws.Range("M135").Value = -3511.3125
$ws.Range("N135").Value = -19323.8814
$ws.Range("H136").Value = 1653.1177
$ws.Range("I136").Value = 1469.2858
$ws.Range("J136").Value = 2511
$ws.Range("K136").Value = 4407.857400000001
$ws.Range("L136").Value = 7533
$ws.Range("M136").Value = 692.1425999999992
$ws.Range("N136").Value = -17733

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1468.5
$ws.Range("I97").Value = 2712
$ws.Range("J97").Value = 777.6667
$ws.Range("K97").Value = 2712
$ws.Range("L97").Value = 777.6667
$ws.Range("M97").Value = -2216
$ws.Range("N97").Value = -1769.6667
$ws.Range("H102").Value = 1665
$ws.Range("I102").Value = 1466.6666
$ws.Range("J102").Value = 1813.75
$ws.Range("K102").Value = 1466.6666
$ws.Range("L102").Value = 1813.75
$ws.Range("M102").Value = 155.3334
$ws.Range("N102").Value = -5057.75

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1510.2106
$ws.Range("I61").Value = 1156.4
$ws.Range("J61").Value = 1903.3334
$ws.Range("K61").Value = 1156.4
$ws.Range("L61").Value = 1903.3334
$ws.Range("M61").Value = -954.4000000000001
$ws.Range("N61").Value = -2307.3334
$ws.Range("H113").Value = 1510.2106
$ws.Range("I113").Value = 1156.4
$ws.Range("J113").Value = 1903.3334
$ws.Range("K113").Value = 1156.4
$ws.Range("L113").Value = 1903.3334
$ws.Range("M113").Value = 1013.6
$ws.Range("N113").Value = -6243.3334
$ws.Range("H132").Value = 22662.77
$ws.Range("I132").Value = 34229.715
$ws.Range("J132").Value = 9168
$ws.Range("K132").Value = 102689.145
$ws.Range("L132").Value = 27504
$ws.Range("M132").Value = -100159.145
$ws.Range("N132").Value = -32564
$ws.Range("H136").Value = 22228314
$ws.Range("I136").Value = 3982.5
$ws.Range("J136").Value = 47627548
$ws.Range("K136").Value = 11947.5
$ws.Range("L136").Value = 142882644
$ws.Range("M136").Value = -9397.5
$ws.Range("N136").Value = -142887744
$ws.Range("H141").Value = 55479.145
$ws.Range("J141").Value = 55479.145
$ws.Range("L141").Value = 55479.145
$ws.Range("N141").Value = -65839.14499999999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1360.4615
$ws.Range("I107").Value = 1192.1
$ws.Range("J107").Value = 1921.6666
$ws.Range("K107").Value = 3576.3
$ws.Range("L107").Value = 5764.9998
$ws.Range("M107").Value = -1656.3
$ws.Range("N107").Value = -9604.9998
$ws.Range("H132").Value = 3733
$ws.Range("I132").Value = 3262.75
$ws.Range("J132").Value = 4270.4287
$ws.Range("K132").Value = 9788.25
$ws.Range("L132").Value = 12811.2861
$ws.Range("M132").Value = -7258.25
$ws.Range("N132").Value = -17871.2861
$ws.Range("H136").Value = 3204.9167
$ws.Range("I136").Value = 1510.3
$ws.Range("J136").Value = 4415.357
$ws.Range("K136").Value = 4530.9
$ws.Range("L136").Value = 13246.071
$ws.Range("M136").Value = -1980.9
$ws.Range("N136").Value = -18346.071

# Remove N132 entirely on ALC sheet (cell deleted per diff)
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("N132").ClearContents()
